$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 13.32779766666667
$ws.Cells.Item(2, 8).Value = 39.983393
$ws.Cells.Item(2, 9).Value = 0.1697233513642653
$ws.Cells.Item(2, 10).Value = 0.1697233513642653
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.116143
$ws.Cells.Item(2, 14).Value = 0.348429
$ws.Cells.Item(2, 15).Value = 0.01430514908838541
$ws.Cells.Item(2, 16).Value = 0.01430514908838541
$ws.Cells.Item(2, 17).Value = 1.547930404399667
$ws.Cells.Item(2, 18).Value = 13.931373639597
$ws.Cells.Item(2, 19).Value = 0.002427917845046236
$ws.Cells.Item(2, 20).Value = 0.002427917845046236

$ws.Cells.Item(3, 7).Value = 13.32779766666667
$ws.Cells.Item(3, 8).Value = 39.983393
$ws.Cells.Item(3, 9).Value = 0.1697233513642653
$ws.Cells.Item(3, 10).Value = 0.1697233513642653
$ws.Cells.Item(3, 15).Value = 0.5605328823946109
$ws.Cells.Item(3, 16).Value = 0.5605328823946107
$ws.Cells.Item(3, 17).Value = 60.65409636512445
$ws.Cells.Item(3, 18).Value = 545.88686728612
$ws.Cells.Item(3, 19).Value = 0.09513551934988494
$ws.Cells.Item(3, 20).Value = 0.09513551934988493

$ws.Cells.Item(4, 7).Value = 13.32779766666667
$ws.Cells.Item(4, 8).Value = 39.983393
$ws.Cells.Item(4, 9).Value = 0.1697233513642653
$ws.Cells.Item(4, 10).Value = 0.1697233513642653
$ws.Cells.Item(4, 15).Value = 0.4251619685170038
$ws.Cells.Item(4, 16).Value = 0.4251619685170038
$ws.Cells.Item(4, 17).Value = 46.005891570625
$ws.Cells.Item(4, 18).Value = 414.053024135625
$ws.Cells.Item(4, 19).Value = 0.07215991416933415
$ws.Cells.Item(4, 20).Value = 0.07215991416933415

$ws.Cells.Item(5, 9).Value = 0.5514955210569645
$ws.Cells.Item(5, 10).Value = 0.5514955210569645
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.116143
$ws.Cells.Item(5, 14).Value = 0.348429
$ws.Cells.Item(5, 15).Value = 0.01430514908838541
$ws.Cells.Item(5, 16).Value = 0.01430514908838541
$ws.Cells.Item(5, 17).Value = 5.029812798724
$ws.Cells.Item(5, 18).Value = 45.26831518851599
$ws.Cells.Item(5, 19).Value = 0.00788922565029667
$ws.Cells.Item(5, 20).Value = 0.00788922565029667

$ws.Cells.Item(6, 9).Value = 0.5514955210569645
$ws.Cells.Item(6, 10).Value = 0.5514955210569645
$ws.Cells.Item(6, 15).Value = 0.5605328823946109
$ws.Cells.Item(6, 16).Value = 0.5605328823946107
$ws.Cells.Item(6, 19).Value = 0.3091313740457781
$ws.Cells.Item(6, 20).Value = 0.3091313740457781

$ws.Cells.Item(7, 9).Value = 0.5514955210569645
$ws.Cells.Item(7, 10).Value = 0.5514955210569645
$ws.Cells.Item(7, 15).Value = 0.4251619685170038
$ws.Cells.Item(7, 16).Value = 0.4251619685170038
$ws.Cells.Item(7, 19).Value = 0.2344749213608898
$ws.Cells.Item(7, 20).Value = 0.2344749213608898

$ws.Cells.Item(8, 8).Value = 65.67520200000001
$ws.Cells.Item(8, 9).Value = 0.27878112757877
$ws.Cells.Item(8, 10).Value = 0.27878112757877
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.116143
$ws.Cells.Item(8, 14).Value = 0.348429
$ws.Cells.Item(8, 15).Value = 0.01430514908838541
$ws.Cells.Item(8, 16).Value = 0.01430514908838541
$ws.Cells.Item(8, 17).Value = 2.542571661962
$ws.Cells.Item(8, 18).Value = 22.883144957658
$ws.Cells.Item(8, 19).Value = 0.003988005593042498
$ws.Cells.Item(8, 20).Value = 0.003988005593042498

$ws.Cells.Item(9, 8).Value = 65.67520200000001
$ws.Cells.Item(9, 9).Value = 0.27878112757877
$ws.Cells.Item(9, 10).Value = 0.27878112757877
$ws.Cells.Item(9, 15).Value = 0.5605328823946109
$ws.Cells.Item(9, 16).Value = 0.5605328823946107
$ws.Cells.Item(9, 17).Value = 99.62811387485337
$ws.Cells.Item(9, 18).Value = 896.6530248736802
$ws.Cells.Item(9, 19).Value = 0.1562659889989477
$ws.Cells.Item(9, 20).Value = 0.1562659889989477

$ws.Cells.Item(10, 8).Value = 65.67520200000001
$ws.Cells.Item(10, 9).Value = 0.27878112757877
$ws.Cells.Item(10, 10).Value = 0.27878112757877
$ws.Cells.Item(10, 15).Value = 0.4251619685170038
$ws.Cells.Item(10, 16).Value = 0.4251619685170038
$ws.Cells.Item(10, 17).Value = 75.56752930125
$ws.Cells.Item(10, 18).Value = 680.1077637112501
$ws.Cells.Item(10, 19).Value = 0.1185271329867799
$ws.Cells.Item(10, 20).Value = 0.1185271329867799
